$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report regeneration for the "647b1ad8-56d5-4878-a597-b65b8878387f"
# file: a new handback (.md) was received but its commit isn't the latest on
# the branch, so the report records the new handback xlf / timestamp and
# raises an "Error Detail" message for both the zh-cn and de-de sheets.
# ---------------------------------------------------------------------------

$mdName   = "647b1ad8-56d5-4878-a597-b65b8878387f.md"
$mdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/13322fef4ec67af651c4873c1d30ab52585f8e65/e2e/647b1ad8-56d5-4878-a597-b65b8878387f.md"
$errorMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae9353cc4615ab8cfc5d75fe96c7ac7dce2588f/e2e/647b1ad8-56d5-4878-a597-b65b8878387f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/13322fef4ec67af651c4873c1d30ab52585f8e65/e2e/647b1ad8-56d5-4878-a597-b65b8878387f.md."

function Update-HandbackRow($ws, $handbackXlf, $handbackDate) {
    # Latest Target File -> link to the (newer) handback markdown file.
    $ws.Range("I6").Value = $mdName
    $ws.Hyperlinks.Add($ws.Range("I6"), $mdUrl, $null, $null, $mdName)
    $ws.Range("I6").Style = "HyperLink"

    # Latest Handback File
    $ws.Range("J6").Value = $handbackXlf

    # Latest Handback DateTime
    $ws.Range("K6").Value = $handbackDate

    # Error Detail
    $ws.Range("P6").Value = $errorMsg

    # Widen the columns that now hold long file names / URLs.
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
    $ws.Columns.Item(16).ColumnWidth = 40
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "647b1ad8-56d5-4878-a597-b65b8878387f.85317950cbeb780e45de6c420ddea764334377cd.zh-cn.xlf" "2016-10-18 03:47:39"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "647b1ad8-56d5-4878-a597-b65b8878387f.85317950cbeb780e45de6c420ddea764334377cd.de-de.xlf" "2016-10-18 03:48:02"
